# Regenerate the "K" column (column G) values on Sheet1.
# The K column previously held "Strike#"-style values; this edit
# recomputes/rewrites those values (s_vals) for every data row (2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 1;  6  = 1;  7  = 0;  8  = 0;  9  = 1;  10 = 0;
    11 = 0;  12 = 0;  13 = 1;  14 = 0;  15 = 0;  16 = 0;  17 = 2;  18 = 0;  19 = 1;
    20 = 0;  21 = 0;  22 = 0;  23 = 1;  24 = 0;  25 = 2;  26 = 1;  27 = 2;  28 = 1;
    29 = 2;  30 = 0;  31 = 2;  32 = 1;  33 = 1;  34 = 1;  35 = 1;  36 = 3;  37 = 0;
    38 = 3;  39 = 1;  40 = 0;  41 = 0;  42 = 1;  43 = 0;  44 = 1;  45 = 1;  46 = 2;
    47 = 0;  48 = 1;  49 = 0;  50 = 2;  51 = 0;  52 = 0;  53 = 2;  54 = 0;  55 = 1;
    56 = 2;  57 = 1;  58 = 2;  59 = 1;  60 = 1;  61 = 0;  62 = 1;  63 = 1;  64 = 1;
    65 = 1;  66 = 1;  67 = 0;  68 = 1;  69 = 2;  70 = 1
}

foreach ($row in 2..70) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
